$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '57.703.11'
Set-TextValue $ws.Range("E2") '  -0.38%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.102.89'
Set-TextValue $ws.Range("E3") '  +1.44%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '523.27'
Set-TextValue $ws.Range("E5") '  +1.29%  '

# Row 6
Set-TextValue $ws.Range("D6") '141.65'
Set-TextValue $ws.Range("E6") '  -0.05%  '

# Row 7
Set-TextValue $ws.Range("E7") '  -0.01%  '

# Row 8
Set-TextValue $ws.Range("D8") '3.103.04'
Set-TextValue $ws.Range("E8") '  +1.52%  '

# Row 9
Set-TextValue $ws.Range("E9") '  +0.46%  '

# Row 10
Set-TextValue $ws.Range("E10") '  -1.51%  '

# Row 11
Set-TextValue $ws.Range("E11") '  +0.23%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.384'
Set-TextValue $ws.Range("E12") '  +2.50%  '

# Row 13
Set-TextValue $ws.Range("E13") '  +1.56%  '

# Row 14
Set-TextValue $ws.Range("D14") '0.132'
Set-TextValue $ws.Range("E14") '  +1.39%  '

# Row 15
Set-TextValue $ws.Range("E15") '  -1.94%  '

# Row 16
Set-TextValue $ws.Range("E16") '  +0.40%  '

# Row 17
Set-TextValue $ws.Range("D17") '57.798.11'
Set-TextValue $ws.Range("E17") '  -0.25%  '

# Row 18
Set-TextValue $ws.Range("D18") '3.109.71'
Set-TextValue $ws.Range("E18") '  +1.91%  '

# Row 19
Set-TextValue $ws.Range("D19") '6.09'
Set-TextValue $ws.Range("E19") '  -0.30%  '

# Row 20
Set-TextValue $ws.Range("D20") '12.76'
Set-TextValue $ws.Range("E20") '  -0.37%  '

# Row 21
Set-TextValue $ws.Range("D21") '8.04'
Set-TextValue $ws.Range("E21") '  -0.21%  '

# Row 22
Set-TextValue $ws.Range("D22") '337.76'
Set-TextValue $ws.Range("E22") '  +1.85%  '

# Row 23
Set-TextValue $ws.Range("E23") '  -0.01%  '

# Row 24
Set-TextValue $ws.Range("E24") '  +1.94%  '

# Row 25
Set-TextValue $ws.Range("D25") '66.52'
Set-TextValue $ws.Range("E25") '  +1.61%  '

# Row 26
Set-TextValue $ws.Range("D26") '0.169'
Set-TextValue $ws.Range("E26") '  -0.62%  '

# Row 27
Set-TextValue $ws.Range("E27") '  +0.21%  '

# Row 28
Set-TextValue $ws.Range("E28") '  +1.40%  '

# Row 29
Set-TextValue $ws.Range("D29") '6.47'
Set-TextValue $ws.Range("E29") '  +0.40%  '

# Row 30
Set-TextValue $ws.Range("E30") '  +0.01%  '

# Row 31
Set-TextValue $ws.Range("D31") '7.16'
Set-TextValue $ws.Range("E31") '  -1.23%  '

# Row 32
Set-TextValue $ws.Range("E32") '  +2.56%  '

# Row 33
Set-TextValue $ws.Range("D33") '20.90'
Set-TextValue $ws.Range("E33") '  +1.06%  '

# Row 34
Set-TextValue $ws.Range("E34") '  -0.97%  '

# Row 35
Set-TextValue $ws.Range("D35") '155.93'
Set-TextValue $ws.Range("E35") '  +0.84%  '

# Row 36
Set-TextValue $ws.Range("E36") '  +1.72%  '

# Row 37
Set-TextValue $ws.Range("D37") '6.11'
Set-TextValue $ws.Range("E37") '  +2.03%  '

# Row 38
Set-TextValue $ws.Range("D38") '27.01'
Set-TextValue $ws.Range("E38") '  +0.34%  '

# Row 39
Set-TextValue $ws.Range("E39") '  -1.18%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.0660'
Set-TextValue $ws.Range("E40") '  -2.83%  '

# Row 41
Set-TextValue $ws.Range("E41") '  +0.76%  '

# Row 42
Set-TextValue $ws.Range("B42") 'RenzoRestakedETH'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue $ws.Range("D42") '3.143.13'
Set-TextValue $ws.Range("E42") '  +1.46%  '

# Row 43
Set-TextValue $ws.Range("B43") 'Stacks'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D43") '1.52'
Set-TextValue $ws.Range("E43") '  +10.89%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.683'
Set-TextValue $ws.Range("E44") '  +3.95%  '

# Row 45
Set-TextValue $ws.Range("D45") '36.77'
Set-TextValue $ws.Range("E45") '  +0.47%  '

# Row 46
Set-TextValue $ws.Range("E46") '  +0.02%  '

# Row 47
Set-TextValue $ws.Range("D47") '2.296.92'
Set-TextValue $ws.Range("E47") '  +0.61%  '

# Row 48
Set-TextValue $ws.Range("E48") '  +0.86%  '

# Row 49
Set-TextValue $ws.Range("E49") '  +4.80%  '

# Row 50
Set-TextValue $ws.Range("D50") '20.50'
Set-TextValue $ws.Range("E50") '  -0.60%  '

# Row 51
Set-TextValue $ws.Range("D51") '6.01'
Set-TextValue $ws.Range("E51") '  +1.50%  '
